$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Date" value (Property/Value pair table) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-01T23:24:47+01:00"

# --- "Include #0" sheet: Concept/Description table ---
$ws = $wb.Worksheets.Item("Include #0")

# Row 2 concept changes from 307157005 / Before breakfast
# to 255214003 / After exercise.
# Assign with a leading apostrophe so the digit-string stays text (not a
# number), then restore the original cell formatting/style (the apostrophe
# trick pulls in an ad-hoc "Text" number format otherwise) by copying the
# format from the row below, which carries the same style as every other
# data row in this table.
$ws.Range("A2").Value = "'255214003"
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2").Value = "After exercise"

# The row for 410594000 / After exercise (row 6) is now redundant and
# gets removed entirely, shifting the remaining rows (7, 8, 9) up by one.
$ws.Rows(6).Delete()
